$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A7").Value = "Sonchus asper"
$ws.Range("B7").Value = "soas"
$ws.Range("E7").Value = "annual"

$ws.Range("H8").Select()
